$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column A: tiny value correction (re-retrieved timestamp)
$ws.Cells.Item(11, 1).Value = 44324.77066125347

# New row 12: next day's data pull
$ws.Cells.Item(12, 1).Value = 44325.77076643925
$ws.Cells.Item(12, 2).Value = 73601
$ws.Cells.Item(12, 3).Value = 61894
$ws.Cells.Item(12, 4).Value = 3275
$ws.Cells.Item(12, 5).Value = 2045
$ws.Cells.Item(12, 6).Value = 1442
$ws.Cells.Item(12, 7).Value = 19160
$ws.Cells.Item(12, 8).Value = 1397
$ws.Cells.Item(12, 9).Value = 835
$ws.Cells.Item(12, 10).Value = 220
